$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "61.487.91"
$ws.Range("E2").Value = "  -1.04%  "

$ws.Range("D3").Value = "3.378.78"
$ws.Range("E3").Value = "  -2.25%  "

$ws.Range("E4").Value = "  +0.09%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "407.10"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.86%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "135.12"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +8.14%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.594"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.04%  "

$ws.Range("E8").Value = "  +0.03%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.672"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +1.82%  "

$ws.Range("E10").Value = "  -4.71%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "42.62"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +2.18%  "

$ws.Range("E12").Value = "  -1.19%  "

$ws.Range("D13").Value = "3.905.81"
$ws.Range("E13").Value = "  -2.13%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "8.41"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -1.55%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "19.72"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -0.61%  "

$ws.Range("D16").Value = "3.361.11"
$ws.Range("E16").Value = "  -2.39%  "

$ws.Range("D17").Value = "61.449.91"
$ws.Range("E17").Value = "  -0.92%  "

$ws.Range("E18").Value = "  -1.60%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "11.00"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.03%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.0000127"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -2.06%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "3.21"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -3.87%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "85.02"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +3.28%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "313.89"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.48%  "

$ws.Range("E24").Value = "  -1.57%  "

$ws.Range("E25").Value = "  -1.10%  "

$ws.Range("E26").Value = "  +11.72%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "8.35"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +5.76%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "29.54"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -5.02%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "7.60"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -3.13%  "

$ws.Range("B30").Value = "Hedera"
$ws.Range("C30").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.117"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +0.24%  "

$ws.Range("B31").Value = "Kaspa"
$ws.Range("C31").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.172"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -1.21%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "2.56"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -0.78%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "11.36"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -2.08%  "

$ws.Range("E34").Value = "  -0.03%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "40.67"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -3.93%  "

$ws.Range("E36").Value = "  -0.73%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "51.88"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -1.05%  "

$ws.Range("E38").Value = "  +0.09%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "3.42"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -2.84%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.93"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -2.30%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "138.61"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +3.15%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.97"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -1.70%  "

$ws.Range("E43").Value = "  -1.10%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.295"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +3.49%  "

$ws.Range("E45").Value = "  +3.21%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "16.77"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -3.11%  "

$ws.Range("E47").Value = "  +1.34%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "21.31"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -4.30%  "

$ws.Range("D49").Value = "2.127.73"
$ws.Range("E49").Value = "  -3.71%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.96"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +3.52%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.28"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -5.41%  "
